$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of pending work item data
$ws.Range("A7").Value = "Bug"
$ws.Range("B7").Value = "Plan Application Main Grid , filters , sort , search not working"
$ws.Range("C7").Value = 45946

# Match the selection left after the edit
$ws.Range("B9").Select()
